$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 133, shifting existing rows 133-177 down to 134-178.
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the weekly record. Most field values
# replicate the former row 133 (now row 134); only the date (D) and volume (M)
# differ for this new entry.
$ws.Cells.Item(133, 1).Value = 11
$ws.Cells.Item(133, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(133, 3).Value = "Bíobío"
$ws.Cells.Item(133, 4).Value = 45120
$ws.Cells.Item(133, 5).Value = 8
$ws.Cells.Item(133, 6).Value = "Fruta"
$ws.Cells.Item(133, 7).Value = 100108
$ws.Cells.Item(133, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(133, 9).Value = 100108002
$ws.Cells.Item(133, 10).Value = "Mango"
$ws.Cells.Item(133, 11).Value = "Sin especificar"
$ws.Cells.Item(133, 12).Value = "Primera"
$ws.Cells.Item(133, 13).Value = 200
$ws.Cells.Item(133, 14).Value = 8000
$ws.Cells.Item(133, 15).Value = 8500
$ws.Cells.Item(133, 16).Value = 8250
$ws.Cells.Item(133, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(133, 18).Value = "Brasil"
$ws.Cells.Item(133, 19).Value = 2062
$ws.Cells.Item(133, 20).Value = 4
